$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old "customer_name" single-cell entry (row 4) ---
# That row is being repurposed as an extra "child" placeholder row, and the
# customer's name/address move down into a new nested "customer" block.
$ws.Range("A4").Value = $null
$ws.Range("G4").Value = $null

# --- Relabel the sample parent.1 / parent.2 values (A/B/C/D -> cell refs) ---
$ws.Range("G2").Value = "G2"
$ws.Range("H2").Value = "H2"
$ws.Range("G3").Value = "G3"
$ws.Range("H3").Value = "H3"

# --- Two additional "child" placeholder rows for json.parent.3 / json.parent.4 ---
$ws.Range("B4").Value = "child"
$ws.Range("B5").Value = "child"
$ws.Range("H5").Value = "H5"

# --- New nested "customer" object (customer.name / customer.address) ---
$ws.Range("A6").Value = "customer"
$ws.Range("B7").Value = "name"
$ws.Range("G7").Value = "山田太郎"
$ws.Range("B8").Value = "address"
$ws.Range("G8").Value = "とうきょう"

# --- Defined names: drop the flat customer_name, add the nested + extra parent names ---
$wb.Names("json.customer_name").Delete()
$wb.Names.Add("json.customer.address", "=Sheet1!`$G`$8")
$wb.Names.Add("json.customer.name", "=Sheet1!`$G`$7")
$wb.Names.Add("json.parent.3", "=Sheet1!`$G`$4:`$H`$4")
$wb.Names.Add("json.parent.4", "=Sheet1!`$G`$5:`$H`$5")

# --- View zoom + print setup ---
$excel.ActiveWindow.Zoom = 115
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
